$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-24 Thursday" "2024-10-25 Friday"

Replace-Text "512×9=" "206×7="
Replace-Text "976×5=" "747×2="
Replace-Text "534×6=" "113×7="
Replace-Text "687×7=" "549×3="
Replace-Text "102×7=" "446×5="

Replace-Text "227×9=" "183×8="
Replace-Text "287×7=" "838×5="
Replace-Text "765×2=" "946×3="
Replace-Text "757×9=" "467×2="
Replace-Text "188×8=" "668×6="

Replace-Text "502×9=" "987×6="
Replace-Text "477×4=" "957×7="
Replace-Text "366×5=" "814×6="
Replace-Text "850×2=" "456×4="
Replace-Text "156×5=" "574×9="

Replace-Text "158×2=" "427×7="
Replace-Text "271×7=" "881×7="
Replace-Text "317×8=" "364×3="
Replace-Text "785×8=" "142×6="
Replace-Text "417×5=" "457×9="

Replace-Text "467×7=" "129×8="
Replace-Text "517×2=" "573×2="
Replace-Text "594×3=" "909×4="
Replace-Text "494×9=" "881×7="
Replace-Text "281×6=" "559×8="
